# Apply updated iteration results to both sheets.
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Full results")
$ws2 = $wb.Worksheets.Item("For plotting")

# --- "Full results" sheet ---
$ws1.Range("H2").Value = 0.567990673360382
$ws1.Range("I2").Value = 0.330304996596306
$ws1.Range("O2").Value = 0.432210294765397

$ws1.Range("F3").Value = 0.569665662819586
$ws1.Range("G3").Value = 0.364314273232692

$ws1.Range("C4").Value = 0.603975554745592
$ws1.Range("D4").Value = 0.396378393439593
$ws1.Range("E4").Value = 1.00035394818519
$ws1.Range("J4").Value = 0.396238145549144
$ws1.Range("K4").Value = 0.364185370454682
$ws1.Range("L4").Value = 0.00167439677955247
$ws1.Range("M4").Value = 0.0359721492162525
$ws1.Range("N4").Value = 0.365859767234235

# --- "For plotting" sheet ---
$ws2.Range("C2").Value = 0.396238145549144
$ws2.Range("D2").Value = 0.337167474042301
$ws2.Range("E2").Value = 0.455308817055987

$ws2.Range("C3").Value = 0.365859767234235
$ws2.Range("D3").Value = 0.302832570431487
$ws2.Range("E3").Value = 0.428886964036983

$ws2.Range("C4").Value = 0.432210294765397
$ws2.Range("D4").Value = 0.369251068965778
$ws2.Range("E4").Value = 0.495169520565016
